$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") cells that are being updated to purely-numeric-looking
# text (e.g. "167.90") must be pre-formatted as Text. Otherwise Excel's
# Range.Value setter auto-converts the string to a real number, which both
# drops significant trailing zeros (e.g. "167.90" -> 167.9) and introduces
# floating point artifacts. NumberFormat is applied per-cell (rather than
# via a single multi-area Range) for reliability.
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).NumberFormat = "@"

$ws.Range("D2").Value = "68.921.35"
$ws.Range("E2").Value = "  +0.92%  "
$ws.Range("D3").Value = "3.746.03"
$ws.Range("E3").Value = "  +1.57%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "602.06"
$ws.Range("E5").Value = "  +0.87%  "
$ws.Range("D6").Value = "167.90"
$ws.Range("E6").Value = "  +0.83%  "
$ws.Range("D7").Value = "3.742.09"
$ws.Range("E7").Value = "  +1.68%  "
$ws.Range("E8").Value = "  -0.10%  "
$ws.Range("E9").Value = "  +1.57%  "
$ws.Range("E10").Value = "  +0.77%  "
$ws.Range("E11").Value = "  +3.29%  "
$ws.Range("E12").Value = "  +0.40%  "
$ws.Range("D13").Value = "38.06"
$ws.Range("E13").Value = "  +0.09%  "
$ws.Range("E14").Value = "  +1.26%  "
$ws.Range("D15").Value = "4.373.19"
$ws.Range("E15").Value = "  +1.49%  "
$ws.Range("D16").Value = "3.746.59"
$ws.Range("E16").Value = "  +1.48%  "
$ws.Range("D17").Value = "68.975.46"
$ws.Range("E17").Value = "  +1.01%  "
$ws.Range("E18").Value = "  +1.40%  "
$ws.Range("E19").Value = "  -1.04%  "
$ws.Range("D20").Value = "17.23"
$ws.Range("E20").Value = "  +0.38%  "
$ws.Range("D21").Value = "10.89"
$ws.Range("E21").Value = "  +19.40%  "
$ws.Range("D22").Value = "493.20"
$ws.Range("E22").Value = "  +0.72%  "
$ws.Range("E23").Value = "  +0.60%  "
$ws.Range("E24").Value = "  +5.90%  "
$ws.Range("D25").Value = "84.85"
$ws.Range("E25").Value = "  +0.49%  "
$ws.Range("D26").Value = "2.31"
$ws.Range("E26").Value = "  +0.82%  "
$ws.Range("D27").Value = "12.30"
$ws.Range("E27").Value = "  +0.83%  "
$ws.Range("E28").Value = "  +0.94%  "
$ws.Range("E29").Value = "  +0.04%  "
$ws.Range("D30").Value = "3.00"
$ws.Range("E30").Value = "  +3.18%  "
$ws.Range("D31").Value = "2.48"
$ws.Range("E31").Value = "  +4.24%  "
$ws.Range("D32").Value = "8.06"
$ws.Range("E32").Value = "  +2.97%  "
$ws.Range("D33").Value = "31.57"
$ws.Range("E33").Value = "  +0.73%  "
$ws.Range("D34").Value = "3.892.01"
$ws.Range("E34").Value = "  +1.55%  "
$ws.Range("E35").Value = "  +0.53%  "
$ws.Range("D36").Value = "3.681.25"
$ws.Range("E36").Value = "  +1.29%  "
$ws.Range("E37").Value = "  -0.02%  "
$ws.Range("D38").Value = "1.02"
$ws.Range("E38").Value = "  +1.89%  "
$ws.Range("E39").Value = "  +2.14%  "
$ws.Range("E40").Value = "  +1.88%  "
$ws.Range("E41").Value = "  +1.18%  "
$ws.Range("D42").Value = "2.96"
$ws.Range("E42").Value = "  +4.95%  "
$ws.Range("D43").Value = "432.48"
$ws.Range("E43").Value = "  -0.25%  "
$ws.Range("D44").Value = "48.60"
$ws.Range("E44").Value = "  -0.02%  "
$ws.Range("E45").Value = "  +2.18%  "
$ws.Range("E46").Value = "  +1.74%  "
$ws.Range("E47").Value = "  -0.01%  "
$ws.Range("D48").Value = "40.47"
$ws.Range("E48").Value = "  +0.27%  "
$ws.Range("D49").Value = "141.33"
$ws.Range("E49").Value = "  -0.11%  "
$ws.Range("D50").Value = "2.790.58"
$ws.Range("E50").Value = "  +1.64%  "
$ws.Range("E51").Value = "  +1.11%  "

# Restore the default ("Normal") cell style on the reformatted cells so the
# only persisted change is the text value itself, matching the original
# (unstyled) Price cells.
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(49, 4).Style = "Normal"

